$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- C3: spMaterialByBatchID note -> renumber last two bullets (3->4->5) ---
# (written before C2 so the shared-string table gets this entry first, matching
#  the author's original commit order)
$c3 = "1. 分析存储过程spMaterialByBatchID：`n  - 尝试在子查询加过滤条件，速度明显提升（30S->5S）；缺陷：在前端不输入过滤条件时，子查询无法过滤数据造成读取次数千万级以上；建议：在不输入查询条件时，仅返回有限的记录数`n  - 尝试在一些连接次数多的表建立索引并改用like xxx%，速度提升40%（5S->2S）；缺陷：只能匹配结尾；`n2. 分析视图vm_BatchStoreQty，尝试减少聚集计算的左外连接数（由4个减为2个）但效果不明显`n3. 压力测试存储过程并监测所占用的进程数和worker数，模拟触发“查询超时已过期”错误`n4. 分析数据连接数，查询等待时间及并行度设置，通过更改配置，调试结果基本无区别。`n5. 压力测试优化过的存储过程，模拟100个客户同时领料查询并根据worker数量控制客户端可以真正进入查询或等待，通过进程的干预暂未触发“查询超时已过期”错误"
$ws.Range("C3").Value = $c3

# --- C2: spsavebaneditlog note -> append advice text, update line 2 wording ---
$c2 = "1. 修改存储过程spsavebaneditlog捕获异常，发现ban_makebill_itemlist_log表字段remotematname长度为50，而ban_makebill_itemlist相应字段长度为500，造成“截断字符串”的错误，在未捕获异常的情况下，前端是否有处理？建议：修改spsavebaneditlog捕获异常并抛出错误信息`n2. 多线程压力测试，100条线程同时调用spsavebaneditlog，在不发生“截断字符串”的错误的情况下，并未发现任何异常。`n3. 归档日志表"
$ws.Range("C2").Value = $c2

# --- C4: sp_pdaOutCaseCommit note -> add advice + new second line, now wraps ---
$c4 = "1. 分析存储过程sp_pdaOutCaseCommit，建议： SET SACT_ABORT ON`n2. 分析“General Network Error,Check your Network Documentation”"
$ws.Range("C4").Value = $c4
$ws.Range("C4").WrapText = $true

# --- Row heights: row 2 grows (more text), row 4 gains an explicit height (now wraps two lines) ---
$ws.Rows.Item(2).RowHeight = 94.5
$ws.Rows.Item(4).RowHeight = 27

# --- Selection moves from C3 to C4 ---
$null = $ws.Range("C4").Select()
